# Update 保險 (insurance) sheet (sheet6) and 事業投資 (business investment)
# sheet (sheet7) with the extended metadata columns (company/name/owner/
# property_category/category/date/legislator_name/legislator_id/
# source_file/index), matching the layout already used on the other
# sheets (土地/建物/汽車/存款/股票).

$wb = $excel.ActiveWorkbook

$legislatorName = "邱議瑩"
$legislatorId = 913
$propertyCategory = "normal"
$reportDate = "2012-04-18"
$sourceFile = "tmped121"

# ---------------------------------------------------------------------
# 保險 sheet (insurance) -- header row + data rows 2-8
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

$insHeader = @("company", "name", "owner", "property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($col = 2; $col -le 11; $col++) {
    $wsIns.Cells.Item(1, $col).Value = $insHeader[$col - 2]
}

$insRows = @(
    @(95, "國泰人壽", "國泰人壽得意還本終身險", "邱議瑩"),
    @(96, "國泰人壽", "國泰人壽富貴年年終身險", "邱議瑩"),
    @(97, "國泰人壽", "國泰人壽添寶養老壽險", "邱議瑩"),
    @(98, "國泰人壽", "國泰人壽雙好還本終身險", "邱議瑩"),
    @(99, "國泰人壽", "國泰人壽雙星還本終身險", "邱議瑩"),
    @(100, "台灣人壽", "台灣人壽歲歲長泰還本终身壽險", "李永得"),
    @(101, "台灣人壽", "台灣人壽新祥和定期壽險", "李永得")
)

$row = 2
foreach ($r in $insRows) {
    $index = $r[0]
    $company = $r[1]
    $name = $r[2]
    $owner = $r[3]

    $wsIns.Cells.Item($row, 1).Value = $index
    $wsIns.Cells.Item($row, 2).Value = $company
    $wsIns.Cells.Item($row, 3).Value = $name
    $wsIns.Cells.Item($row, 4).Value = $owner
    $wsIns.Cells.Item($row, 5).Value = "insurance"
    $wsIns.Cells.Item($row, 6).Value = $propertyCategory
    $wsIns.Cells.Item($row, 7).Value = $reportDate
    $wsIns.Cells.Item($row, 8).Value = $legislatorName
    $wsIns.Cells.Item($row, 9).Value = $legislatorId
    $wsIns.Cells.Item($row, 10).Value = $sourceFile
    $wsIns.Cells.Item($row, 11).Value = $index

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 事業投資 sheet (business investment) -- header row + data rows 2-4
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("事業投資")

$invHeader = @("owner", "company", "address", "total", "register_date", "register_reason", "property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($col = 2; $col -le 14; $col++) {
    $wsInv.Cells.Item(1, $col).Value = $invHeader[$col - 2]
}

$invRows = @(
    @(116, "李永得", "玉山社事業股份有限公司", "臺北市仁愛路四段145號3樓之2", 1000000, "84年07月08日", "發起設立"),
    @(117, "李永得", "淡海股份有限公司", "新北市真理街17號", 280800, "97年02月27日", "發起設立"),
    @(118, "李永得", "先驅媒體社會企業股份有限公司", "臺北市仁愛路二段98號7樓", 500000, "98年07月24日", "發起設立")
)

$row = 2
foreach ($r in $invRows) {
    $index = $r[0]
    $owner = $r[1]
    $company = $r[2]
    $address = $r[3]
    $total = $r[4]
    $registerDate = $r[5]
    $registerReason = $r[6]

    $wsInv.Cells.Item($row, 1).Value = $index
    $wsInv.Cells.Item($row, 2).Value = $owner
    $wsInv.Cells.Item($row, 3).Value = $company
    $wsInv.Cells.Item($row, 4).Value = $address
    $wsInv.Cells.Item($row, 5).Value = $total
    $wsInv.Cells.Item($row, 6).Value = $registerDate
    $wsInv.Cells.Item($row, 7).Value = $registerReason
    $wsInv.Cells.Item($row, 8).Value = "investment"
    $wsInv.Cells.Item($row, 9).Value = $propertyCategory
    $wsInv.Cells.Item($row, 10).Value = $reportDate
    $wsInv.Cells.Item($row, 11).Value = $legislatorName
    $wsInv.Cells.Item($row, 12).Value = $legislatorId
    $wsInv.Cells.Item($row, 13).Value = $sourceFile
    $wsInv.Cells.Item($row, 14).Value = $index

    $row = $row + 1
}
